$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "95.871.39"
$ws.Range("E2").Value = "  -0.75%  "

# Row 3
$ws.Range("D3").Value = "3.469.30"
$ws.Range("E3").Value = "  +4.27%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.03"
$ws.Range("E5").Value = "  -2.28%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "645.94"
$ws.Range("E6").Value = "  -0.83%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.46"
$ws.Range("E7").Value = "  +7.72%  "

# Row 8
$ws.Range("E8").Value = "  -0.76%  "

# Row 9
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.995"
$ws.Range("E10").Value = "  +1.63%  "

# Row 11
$ws.Range("D11").Value = "3.469.40"
$ws.Range("E11").Value = "  +4.33%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.74"
$ws.Range("E12").Value = "  +7.80%  "

# Row 13
$ws.Range("E13").Value = "  -3.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.14"
$ws.Range("E14").Value = "  -0.35%  "

# Row 15
$ws.Range("D15").Value = "95.692.74"
$ws.Range("E15").Value = "  -0.67%  "

# Row 16
$ws.Range("D16").Value = "4.125.82"
$ws.Range("E16").Value = "  +4.45%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000255"
$ws.Range("E17").Value = "  +1.89%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.46"
$ws.Range("E18").Value = "  +0.08%  "

# Row 19
$ws.Range("D19").Value = "3.469.83"
$ws.Range("E19").Value = "  +4.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.27"
$ws.Range("E20").Value = "  +9.32%  "

# Row 21
$ws.Range("E21").Value = "  +13.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.516"
$ws.Range("E22").Value = "  +7.78%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "510.15"
$ws.Range("E23").Value = "  +0.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.24"
$ws.Range("E24").Value = "  -2.52%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000193"
$ws.Range("E25").Value = "  -1.68%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.65"
$ws.Range("E26").Value = "  +3.55%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "92.27"
$ws.Range("E27").Value = "  -2.29%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.32"
$ws.Range("E28").Value = "  +3.44%  "

# Row 29
$ws.Range("D29").Value = "3.661.74"
$ws.Range("E29").Value = "  +4.59%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.83"
$ws.Range("E30").Value = "  +10.34%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.28%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.76"
$ws.Range("E32").Value = "  +12.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.138"
$ws.Range("E33").Value = "  -2.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.184"
$ws.Range("E34").Value = "  -0.57%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "31.00"
$ws.Range("E35").Value = "  +11.17%  "

# Row 36
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.10%  "

# Row 37
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.575"
$ws.Range("E37").Value = "  +6.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.82"
$ws.Range("E38").Value = "  +4.32%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.46"
$ws.Range("E39").Value = "  -2.84%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "517.87"
$ws.Range("E40").Value = "  +2.47%  "

# Row 41
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.12%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.150"
$ws.Range("E42").Value = "  +0.40%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.915"
$ws.Range("E43").Value = "  +11.27%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.14"
$ws.Range("E44").Value = "  -1.41%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.71"
$ws.Range("E45").Value = "  +6.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0419"
$ws.Range("E46").Value = "  +3.69%  "

# Row 47
$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.62"
$ws.Range("E47").Value = "  -0.83%  "

# Row 48
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.57"
$ws.Range("E48").Value = "  +3.25%  "

# Row 49
$ws.Range("E49").Value = "  +11.69%  "

# Row 50
$ws.Range("E50").Value = "  +3.23%  "

# Row 51
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.62"
$ws.Range("E51").Value = "  +1.11%  "
